# Added get_form_type in the import_utils
#
# Adds a new "Form Tag" column (column V) to the capital-commitments
# allocation sheet and defaults every existing data row (2-8) to "Default".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("V1").Value = "Form Tag"

# Default tag value for each existing data row
$ws.Range("V2:V8").Value = "Default"

# Leave the selection on the newly added column, like the author did
$ws.Range("V2:V8").Select()
